$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 100: add a status note in column D
$ws.Range("D100").Value = "done 1 error in logical or and"

# Row 101: replace old "coding" value (moved down to row 107) with new entries
$ws.Range("E101").Value = "x - 3:10"
$ws.Range("D101").Value = "coded and tested"

# Row 102: new topic entry with hyperlink in column A
$ws.Range("C102").Value = "How to turn off a particular bit in a number?"
$ws.Hyperlinks.Add($ws.Range("A102"), "https://www.geeksforgeeks.org/how-to-turn-off-a-particular-bit-in-a-number/")
$ws.Range("A102").Style = "Hyperlink"
$ws.Range("E102").Value = "3:22 - 3:37"
$ws.Range("D102").Value = "done"

# Row 103
$ws.Range("E103").Value = "x - 3:42"
$ws.Range("D103").Value = "coded"

# Row 104
$ws.Range("C104").Value = "Find Excel column name from a given column number"
$ws.Range("E104").Value = "3:56 - 4:13"
$ws.Range("D104").Value = "skipped"

# Row 105: new topic entry with hyperlink in column A
$ws.Range("C105").Value = "Program to print all palindromes in a given range"
$ws.Hyperlinks.Add($ws.Range("A105"), "https://www.geeksforgeeks.org/program-print-palindromes-given-range/")
$ws.Range("A105").Style = "Hyperlink"
$ws.Range("E105").Value = "4:19 - 4:34"
$ws.Range("D105").Value = "done brute force"

# Row 106: new topic entry with hyperlink in column A
$ws.Range("C106").Value = "Print all pairs of anagrams in a given array of strings"
$ws.Hyperlinks.Add($ws.Range("A106"), "https://www.geeksforgeeks.org/print-pairs-anagrams-given-array-strings/")
$ws.Range("A106").Style = "Hyperlink"
$ws.Range("E106").Value = "4:49 - 5:22"

# Row 107: the old "coding" note that used to live in D101
$ws.Range("D107").Value = "coding"

# Update view to match new extent
$ws.Range("E107").Select()
$excel.ActiveWindow.ScrollRow = 94
